$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update due-date values (column C) that were edited in the source workbook
$ws.Range("C2").Value = 45296
$ws.Range("C3").Value = 45475
$ws.Range("C4").Value = 45471
$ws.Range("C7").Value = 45506
$ws.Range("C9").Value = 45319
$ws.Range("C10").Value = 45293

# Helper column E: whether the due date (column C) is already in the past
$ws.Range("E2").Formula = "=C2<TODAY()"
$ws.Range("E3:E10").Formula = "=C3<TODAY()"

# Conditional formatting: highlight overdue due dates in C2:C10 with a bold
# white-on-red fill, driven by the same expression as column E.
$cf = $ws.Range("C2:C10").FormatConditions.Add(2, 0, "=C2<TODAY()")
$cf.Font.Bold = $true
$cf.Font.Color = 16777215
$cf.Interior.Color = 255

# Restore the view state recorded in the edited workbook
[void]$ws.Range("D17").Select()
$excel.ActiveWindow.Zoom = 115
